# "Merge of harris into stable"
#
# The manifest sheet had a duplicate shared string
# ("handbook_of_archaeology.jpg") that the real file's header row
# ("filename") should have matched to. Re-casing the stray duplicate so
# it collapses into a single shared-string entry, widening the filename
# column so it is readable, and moving the selection/cursor down past the
# header+data rows to A4 (where a user would continue entering data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-cased filename value in row 3 (the duplicate "handbook_..."
# shared string collapses into the corrected "HANDBOOK_..." one on save).
$ws.Range("A3").Value = "HANDBOOK_of_archaeology.jpg"

# Widen column A (filenames are long) and leave it as an explicit custom
# width (~34.57 chars).
$ws.Columns.Item(1).ColumnWidth = 33.7

# Move the active selection down to the first empty row below the data.
$ws.Range("A4").Select()
